{"js": "// Find the paragraph describing the \"\u0418\u0441\u0442\u043e\u0440\u0438\u044f \u0437\u0430\u043a\u0430\u0437\u043e\u0432\" (order history) section\n// for the regular user, and append \", \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\" right after the existing\n// \"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\" status list, so the sentence reads:\n// \"... \u0441\u0442\u0430\u0442\u0443\u0441 \u0437\u0430\u043a\u0430\u0437\u0430 (\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb), \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\"\nconst body = context.document.body;\n\nconst results = body.search(\"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Anchor text \"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\" not found in document body.');\n}\n\n// Insert the new text immediately after the matched anchor text.\nconst anchor = results.items[0];\nanchor.insertText(\", \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Find the paragraph describing the \"\u0418\u0441\u0442\u043e\u0440\u0438\u044f \u0437\u0430\u043a\u0430\u0437\u043e\u0432\" (order history) section\n# for the regular user, and append \", \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\" right after the existing\n# \"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\" status list, so the sentence reads:\n# \"... \u0441\u0442\u0430\u0442\u0443\u0441 \u0437\u0430\u043a\u0430\u0437\u0430 (\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb), \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\"\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindStop\n\n$found = $find.Execute()\n\nif ($found) {\n    $range.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n    $range.InsertAfter(\", \u0434\u0430\u0442\u0443 \u0437\u0430\u043a\u0430\u0437\u0430\")\n} else {\n    throw 'Anchor text \"(\u00ab\u041f\u043e\u043b\u0443\u0447\u0435\u043d\u00bb, \u00ab\u041e\u0442\u043c\u0435\u043d\u0435\u043d\u00bb)\" not found in document.'\n}\n"}
